# Apply cryptos list update (values scraped on Wed May  8 21:46:23 UTC 2024)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = "Normal"
}

Set-TextValue "D2" "61.654.91"
Set-TextValue "E2" "  -2.15%  "
Set-TextValue "D3" "2.978.60"
Set-TextValue "E3" "  -2.36%  "
Set-TextValue "E4" "  -0.03%  "
Set-TextValue "D5" "588.29"
Set-TextValue "E5" "  +1.00%  "
Set-TextValue "D6" "141.97"
Set-TextValue "E6" "  -5.98%  "
Set-TextValue "E7" "  -0.02%  "
Set-TextValue "E8" "  -2.64%  "
Set-TextValue "D9" "2.976.67"
Set-TextValue "E9" "  -2.39%  "
Set-TextValue "D10" "0.143"
Set-TextValue "E10" "  -6.02%  "
Set-TextValue "D11" "5.77"
Set-TextValue "E11" "  -0.88%  "
Set-TextValue "E12" "  +2.42%  "
Set-TextValue "E13" "  -3.10%  "
Set-TextValue "D14" "34.01"
Set-TextValue "E14" "  -5.63%  "
Set-TextValue "E15" "  +1.33%  "
Set-TextValue "D16" "3.464.99"
Set-TextValue "E16" "  -2.44%  "
Set-TextValue "D17" "7.01"
Set-TextValue "E17" "  -1.70%  "
Set-TextValue "D18" "61.592.47"
Set-TextValue "E18" "  -2.19%  "
Set-TextValue "D19" "2.973.63"
Set-TextValue "E19" "  -2.58%  "
Set-TextValue "D20" "450.16"
Set-TextValue "E20" "  -6.25%  "
Set-TextValue "D21" "13.88"
Set-TextValue "E21" "  -2.89%  "
Set-TextValue "D22" "0.683"
Set-TextValue "E22" "  -3.26%  "
Set-TextValue "D23" "7.32"
Set-TextValue "E23" "  -2.66%  "
Set-TextValue "D24" "81.16"
Set-TextValue "E24" "  -0.86%  "
Set-TextValue "D25" "12.12"
Set-TextValue "E25" "  -4.01%  "
Set-TextValue "E26" "  -9.89%  "
Set-TextValue "E27" "  +0.10%  "
Set-TextValue "D28" "9.82"
Set-TextValue "E28" "  -6.86%  "
Set-TextValue "E29" "  -0.08%  "
Set-TextValue "E30" "  -0.60%  "
Set-TextValue "E31" "  -7.11%  "
Set-TextValue "E32" "  -6.10%  "
Set-TextValue "E33" "  -2.12%  "
Set-TextValue "E34" "  -3.25%  "
Set-TextValue "E35" "  -4.86%  "
Set-TextValue "D36" "0.0₃0775"
Set-TextValue "E36" "  -4.43%  "
Set-TextValue "D37" "5.69"
Set-TextValue "E37" "  -3.72%  "
Set-TextValue "E38" "  -5.11%  "
Set-TextValue "B39" "Cosmos"
Set-TextValue "C39" "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
Set-TextValue "D39" "9.14"
Set-TextValue "E39" "  -0.29%  "
Set-TextValue "B40" "OKB"
Set-TextValue "C40" "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
Set-TextValue "D40" "50.12"
Set-TextValue "E40" "  -0.43%  "
Set-TextValue "B41" "dogwifhat"
Set-TextValue "C41" "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
Set-TextValue "D41" "2.81"
Set-TextValue "E41" "  -10.37%  "
Set-TextValue "B42" "Kaspa"
Set-TextValue "C42" "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
Set-TextValue "D42" "0.119"
Set-TextValue "E42" "  +3.11%  "
Set-TextValue "D43" "387.75"
Set-TextValue "E43" "  -8.87%  "
Set-TextValue "D44" "0.0354"
Set-TextValue "E44" "  -1.96%  "
Set-TextValue "D45" "2.724.90"
Set-TextValue "E45" "  -4.23%  "
Set-TextValue "E46" "  -8.50%  "
Set-TextValue "D47" "36.93"
Set-TextValue "E47" "  -2.20%  "
Set-TextValue "D48" "130.12"
Set-TextValue "E48" "  +2.53%  "
Set-TextValue "E49" "  +0.08%  "
Set-TextValue "D50" "0.108"
Set-TextValue "E50" "  -1.71%  "
Set-TextValue "E51" "  -1.03%  "
